$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Row 7: add Description for InputFolder (wrap-text style like C3/C5) ---
$ws.Range("C7").Value = "Input Directory for invoice files"
$ws.Range("C7").WrapText = $true

# --- Row 8: add Description for OutputFile ---
$ws.Range("C8").Value = "output file for summary"

# --- Row 9: new setting MispellingTolerance ---
$ws.Range("A9").Value = "MispellingTolerance"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "maximum character misspellings allowed in a field name"
$ws.Range("C9").WrapText = $true

# --- Row 10: new setting SortField ---
$ws.Range("A10").Value = "SortField"
$ws.Range("B10").Value = "Due Date"
$ws.Range("C10").Value = "field by which to sort datatable output"

# Row 11 stays blank

# --- Row 12: new setting APIKey ---
$ws.Range("A12").Value = "APIKey"
$ws.Range("B12").Value = "3uArrqaAfpX1h7h5JQasmNxmuJ5eTekx"
$ws.Range("C12").Value = "API Key for currency converter"

# --- Row 13: new setting Endpoint ---
$ws.Range("A13").Value = "Endpoint"
$ws.Range("B13").Value = "https://api.apilayer.com/exchangerates_data/convert"
$ws.Range("C13").Value = "URL for Conversion"

# --- Update the active selection to A13 (last edited row) ---
$ws.Activate()
$ws.Range("A13").Select()
